$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 375 (shifts old rows 375-422 down to 376-423)
$ws.Rows.Item(375).Insert()

# Populate the newly inserted row 375 with the new weekly data point
$ws.Cells.Item(375, 1).Value  = 5
$ws.Cells.Item(375, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(375, 3).Value  = "Maule"
$ws.Cells.Item(375, 4).Value  = 44951
$ws.Cells.Item(375, 5).Value  = 7
$ws.Cells.Item(375, 6).Value  = 100112006
$ws.Cells.Item(375, 7).Value  = "Repollo"
$ws.Cells.Item(375, 8).Value  = "Crespo record"
$ws.Cells.Item(375, 9).Value  = "Primera"
$ws.Cells.Item(375, 10).Value = 1500
$ws.Cells.Item(375, 11).Value = 1200
$ws.Cells.Item(375, 12).Value = 1200
$ws.Cells.Item(375, 13).Value = 1200
$ws.Cells.Item(375, 14).Value = "$/unidad"
$ws.Cells.Item(375, 15).Value = "Región del Maule"
$ws.Cells.Item(375, 16).Value = 1200
$ws.Cells.Item(375, 17).Value = 1
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# Make sure the date cell carries the same number format as the rest of column D (s="2")
$ws.Cells.Item(375, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
